# Tutorial 9 slides fix
# Slide 10, "Content Placeholder 2": the second "Binary Search" in
#   "Binary Search is much more efficient than Binary Search."
# was actually meant to read "Linear Search" (the bullet is contrasting
# Binary Search against Linear Search), e.g.:
#   "Binary Search is much more efficient than Linear Search."

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)

# Locate the content placeholder shape that holds the "Binary Search ..."
# bullet text rather than relying on a hard-coded shape index.
$targetShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTextFrame -and $candidate.TextFrame.HasText) {
        if ($candidate.TextFrame.TextRange.Text -like "*Binary Search*") {
            $targetShape = $candidate
        }
    }
}

$tr = $targetShape.TextFrame.TextRange
$para1 = $tr.Paragraphs(1)
$paraText = $para1.Text

$searchTerm = "Binary Search"
$lastIdx = $paraText.LastIndexOf($searchTerm)
$absoluteStart = $para1.Start + $lastIdx

$targetRun = $tr.Characters($absoluteStart, $searchTerm.Length)
$targetRun.Text = "Linear Search"
